$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Cells.Item(2, 4).Style
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '67.382.12'
$ws.Cells.Item(2, 4).Style = $origStyle
$ws.Cells.Item(2, 5).Value = '  +0.48%  '

$origStyle = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.634.52'
$ws.Cells.Item(3, 4).Style = $origStyle
$ws.Cells.Item(3, 5).Value = '  +0.73%  '

$ws.Cells.Item(4, 5).Value = '  +0.15%  '

$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '602.14'
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  +1.62%  '

$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '153.37'
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  -0.45%  '

$ws.Cells.Item(7, 5).Value = '  +0.03%  '

$origStyle = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.560'
$ws.Cells.Item(8, 4).Style = $origStyle
$ws.Cells.Item(8, 5).Value = '  +3.60%  '

$origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.633.01'
$ws.Cells.Item(9, 4).Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  +0.67%  '

$ws.Cells.Item(10, 5).Value = '  +5.49%  '

$ws.Cells.Item(11, 5).Value = '  +0.50%  '

$origStyle = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.21'
$ws.Cells.Item(12, 4).Style = $origStyle
$ws.Cells.Item(12, 5).Value = '  -0.02%  '

$origStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.353'
$ws.Cells.Item(13, 4).Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  -0.70%  '

$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '27.83'
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  +0.34%  '

$origStyle = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.113.69'
$ws.Cells.Item(15, 4).Style = $origStyle
$ws.Cells.Item(15, 5).Value = '  +0.82%  '

$ws.Cells.Item(16, 5).Value = '  +0.42%  '

$origStyle = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '67.273.30'
$ws.Cells.Item(17, 4).Style = $origStyle
$ws.Cells.Item(17, 5).Value = '  +0.24%  '

$origStyle = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.638.83'
$ws.Cells.Item(18, 4).Style = $origStyle
$ws.Cells.Item(18, 5).Value = '  +1.05%  '

$origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.25'
$ws.Cells.Item(19, 4).Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  -0.25%  '

$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '364.17'
$ws.Cells.Item(20, 4).Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  +1.38%  '

$ws.Cells.Item(21, 5).Value = '  -3.91%  '

$ws.Cells.Item(22, 5).Value = '  -0.41%  '

$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '2.14'
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  +4.63%  '

$ws.Cells.Item(24, 5).Value = '  +0.05%  '

$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '10.16'
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  -1.02%  '

$origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '67.38'
$ws.Cells.Item(26, 4).Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  -5.68%  '

$ws.Cells.Item(27, 2).Value = 'WrappedeETH'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.763.94'
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  +0.54%  '

$ws.Cells.Item(28, 2).Value = 'PEPE'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$origStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.0000104'
$ws.Cells.Item(28, 4).Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  +0.06%  '

$ws.Cells.Item(29, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$origStyle = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 4).Style = $origStyle
$ws.Cells.Item(29, 5).Value = '  +0.18%  '

$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '576.59'
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  -7.74%  '

$ws.Cells.Item(31, 5).Value = '  -3.92%  '

$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.88'
$ws.Cells.Item(32, 4).Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  -1.55%  '

$ws.Cells.Item(33, 5).Value = '  -0.19%  '

$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.129'
$ws.Cells.Item(34, 4).Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  -3.80%  '

$ws.Cells.Item(35, 5).Value = '  +0.10%  '

$ws.Cells.Item(36, 5).Value = '  -2.18%  '

$ws.Cells.Item(37, 5).Value = '  -1.29%  '

$origStyle = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '157.64'
$ws.Cells.Item(38, 4).Style = $origStyle
$ws.Cells.Item(38, 5).Value = '  +2.56%  '

$origStyle = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '19.42'
$ws.Cells.Item(39, 4).Style = $origStyle

$ws.Cells.Item(40, 5).Value = '  -0.04%  '

$origStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.29'
$ws.Cells.Item(41, 4).Style = $origStyle
$ws.Cells.Item(41, 5).Value = '  -3.91%  '

$ws.Cells.Item(42, 5).Value = '  -0.47%  '

$origStyle = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.64'
$ws.Cells.Item(43, 4).Style = $origStyle
$ws.Cells.Item(43, 5).Value = '  +1.37%  '

$origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '41.23'
$ws.Cells.Item(44, 4).Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  -0.22%  '

$ws.Cells.Item(45, 5).Value = '  +0.01%  '

$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '16.36'
$ws.Cells.Item(46, 4).Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  -0.80%  '

$origStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '156.26'
$ws.Cells.Item(47, 4).Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  +0.38%  '

$ws.Cells.Item(48, 5).Value = '  -3.62%  '

$origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '3.74'
$ws.Cells.Item(49, 4).Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  -1.03%  '

$ws.Cells.Item(50, 2).Value = 'Mantle'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.625'
$ws.Cells.Item(50, 4).Style = $origStyle
$ws.Cells.Item(50, 5).Value = '  +0.21%  '

$ws.Cells.Item(51, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '20.88'
$ws.Cells.Item(51, 4).Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  -0.80%  '
